# Chapter 4 - working... (word_list.xlsx update)
# Adds a couple of new glossary terms and fills in a previously partial row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 93: supremum / 상계, flagged for follow-up in the 메모 column.
$ws.Range("A93").Value = "supremum"
$ws.Range("B93").Value = "상계"
$ws.Range("D93").Value = "확인필요"

# Row 62 (power series / 멱급수) gains entries in the 번역2 (C) and 메모 (D) columns.
$ws.Range("C62").Value = "제곱급수"
$ws.Range("D62").Value = "제곱급수"

# New row 94: complex power series / 복소제곱급수, with a spacing question noted.
$ws.Range("A94").Value = "complex power series"
$ws.Range("B94").Value = "복소제곱급수"
$ws.Range("D94").Value = "띄어쓰기는?"

# Leave the cursor where the author left off: the empty cell below the new rows.
$ws.Range("B95").Select() | Out-Null
